$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)
$ws3 = $wb.Worksheets.Item(3)
$ws5 = $wb.Worksheets.Item(5)

# =========================================================
# OrderProducts table (rows 19-26)
# =========================================================

# Title cell (bold font, no fill) - copy format from an existing "table" title cell
$ws3.Range("A2").Copy()
$ws4.Range("A19").PasteSpecial(-4122)
$ws4.Range("A19").Value = "OrderProducts table"
$ws4.Range("A19").Characters(14, 6).Font.Bold = $false

# Header row (bold + fill)
$ws4.Range("A3:C3").Copy()
$ws4.Range("A20:C20").PasteSpecial(-4122)
$ws4.Range("A20").Value = "OrderID"
$ws4.Range("B20").Value = "ProductID"
$ws4.Range("C20").Value = "Quantity"

# Data rows
$ws4.Range("A21").Value = 1
$ws4.Range("B21").Value = "CB-2903"
$ws4.Range("C21").Value = 1

$ws4.Range("A22").Value = 1
$ws4.Range("B22").Value = "BA-3827"
$ws4.Range("C22").Value = 2

$ws4.Range("A23").Value = 2
$ws4.Range("B23").Value = "BA-3827"
$ws4.Range("C23").Value = 1

$ws4.Range("A24").Value = 2
$ws4.Range("B24").Value = "BA-2349"
$ws4.Range("C24").Value = 1

$ws4.Range("A25").Value = 2
$ws4.Range("B25").Value = "BA-2908"
$ws4.Range("C25").Value = 1

$ws4.Range("A26").Value = 3
$ws4.Range("B26").Value = "BE-2349"
$ws4.Range("C26").Value = 1

# =========================================================
# Orders table (row 12, cols A-C) and Products table (row 12, cols E-F)
# =========================================================

# Title cells
$ws4.Range("A2").Copy()
$ws4.Range("A12").PasteSpecial(-4122)
$ws4.Range("A12").Value = "Orders table"

$ws3.Range("A14").Copy()
$ws4.Range("E12").PasteSpecial(-4122)
$ws4.Range("E12").Value = "Products table"
$ws4.Range("E12").Characters(9, 6).Font.Bold = $false

# Header row 13
$ws4.Range("A3:C3").Copy()
$ws4.Range("A13:C13").PasteSpecial(-4122)
$ws4.Range("A13").Value = "OrderID"
$ws4.Range("B13").Value = "OrderDate"
$ws4.Range("C13").Value = "CustomerID"

$ws4.Range("A3:B3").Copy()
$ws4.Range("E13:F13").PasteSpecial(-4122)
$ws4.Range("E13").Value = "ProductID"
$ws4.Range("F13").Value = "Price"

# Date formatting for B14:B16 (copy from an existing OrderDate column)
$ws5.Range("B4").Copy()
$ws4.Range("B14:B16").PasteSpecial(-4122)

# Orders data
$ws4.Range("A14").Value = 1
$ws4.Range("B14").Value = 42371
$ws4.Range("C14").Value = 101

$ws4.Range("A15").Value = 2
$ws4.Range("B15").Value = 42371
$ws4.Range("C15").Value = 163

$ws4.Range("A16").Value = 3
$ws4.Range("B16").Value = 42372
$ws4.Range("C16").Value = 302

# Products data
$ws4.Range("E14").Value = "CB-2903"
$ws4.Range("F14").Value = 12.99
$ws4.Range("F14").NumberFormat = "0.00"

$ws4.Range("E15").Value = "BA-3827"
$ws4.Range("F15").Value = 1.5

$ws4.Range("E16").Value = "BA-2349"
$ws4.Range("F16").Value = 5.99

$ws4.Range("E17").Value = "BA-2903"
$ws4.Range("F17").Value = 10

# Apply the new price number format (0.00) to F15:F17 to match F14
$ws4.Range("F14").Copy()
$ws4.Range("F15:F17").PasteSpecial(-4122)
$ws4.Range("F15").Value = 1.5
$ws4.Range("F16").Value = 5.99
$ws4.Range("F17").Value = 10

# =========================================================
# Sheet selection / active tab
# =========================================================
$ws4.Range("H13").Select()
$ws4.Activate()

# =========================================================
# Page setup
# =========================================================
$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1
